$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows by copy-inserting row 22, which carries over the
# existing cell styles (date format on column B, numeric style on column C)
# onto the freshly created rows 23 and 24.
$ws.Rows("22:22").Copy()
$ws.Rows("23:23").Insert(-4121)
$ws.Rows("22:22").Copy()
$ws.Rows("24:24").Insert(-4121)
$excel.CutCopyMode = 0

# New row 23: A=21, B=23/06/2025 (date serial 45831), C=642, D=5, E="Ronaldo"
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 45831
$ws.Cells.Item(23, 3).Value = 642
$ws.Cells.Item(23, 4).Value = 5
$ws.Cells.Item(23, 5).Value = "Ronaldo"

# New row 24: A=22, B=23/06/2025 (date serial 45831), C=642, D=5, E="Bruno"
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 45831
$ws.Cells.Item(24, 3).Value = 642
$ws.Cells.Item(24, 4).Value = 5
$ws.Cells.Item(24, 5).Value = "Bruno"

# Update selection to match the final cursor position
$ws.Range("E24").Select()
